$wb = $excel.ActiveWorkbook

# Sheet: Folder Inventory - C2 timestamp update
$wsFolderInventory = $wb.Worksheets.Item("Folder Inventory")
$wsFolderInventory.Range("C2").Value = "2025-06-12 17:37:08 +0530"

# Sheet: Metadata - B3 "Generated On" timestamp update, B5 "Workflow Run" counter update
$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsMetadata.Range("B3").Value = "2025-06-12 12:07:25 UTC"
# B5 ("Workflow Run") is stored as text (e.g. "14"), not a number, so force
# the cell to Text format first - otherwise a digit-only string like "15"
# would be auto-converted to a numeric value.
$wsMetadata.Range("B5").NumberFormat = "@"
$wsMetadata.Range("B5").Value = "15"

# Sheet: Summary - B5 "Most Recent Update" timestamp update
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B5").Value = "2025-06-12 17:37:08 +0530"
